$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) "UNO Yocto+SPDX Group -Tom Gurney, Tim Strever, Dan Right, and Kevin
#    Lumbard" -> split into three runs and correct "Right" to "Wright":
#       "UNO Yocto+SPDX Group"
#       " -Tom Gurney, Tim Strever, Dan Wr"
#       "ight, and Kevin Lumbard"
# -----------------------------------------------------------------------
$whole = $d.Content
$whole.Find.Execute("UNO Yocto+SPDX Group -Tom Gurney, Tim Strever, Dan Right, and Kevin Lumbard", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$wholeStart = $whole.Start

# Fix the typo: replace the "R" in "Right" with "Wr" (Right -> Wright)
$rPos = $wholeStart + 51
$rCharRange = $d.Range($rPos, $rPos + 1)
$rCharRange.InsertBefore("Wr")
$rDel = $d.Range($rPos + 2, $rPos + 3)
$rDel.Delete()

# Force the run to split at the two boundaries by nudging (and restoring)
# direct formatting on exactly the middle segment - this creates distinct
# runs without altering the visible formatting.
$midStart = $wholeStart + 20
$midEnd = $rPos + 2
$mid = $d.Range($midStart, $midEnd)
$mid.Bold = 1
$mid.Bold = 0

# -----------------------------------------------------------------------
# 2) "Day to Da" + hidden _GoBack bookmark + "y Business"
#    -> single run "Day to Day Business" (bookmark removed from here)
# -----------------------------------------------------------------------
$dayRng = $d.Content
$dayRng.Find.Execute("Day to Da", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dayStart = $dayRng.Start

# Two-step text replacement: first to a placeholder (forces the engine to
# actually commit the structural merge/bookmark removal), then to the
# final text.
$fullDayRange = $d.Range($dayStart, $dayStart + 19)
$fullDayRange.Text = "xxxxxxxxxxxxxxxxxxx"
$fullDayRange2 = $d.Range($dayStart, $dayStart + 19)
$fullDayRange2.Text = "Day to Day Business"

# -----------------------------------------------------------------------
# 3) Re-create the _GoBack bookmark right after " Yocto - "
# -----------------------------------------------------------------------
$yocto = $d.Content
$yocto.Find.Execute(" Yocto - ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $yocto.End

# Bookmarks.Add mis-handles a truly collapsed Range, so insert a one
# character placeholder, anchor the bookmark to it, then delete the
# placeholder - leaving a collapsed bookmark at the right spot.
$insRange = $d.Range($insertPos, $insertPos)
$insRange.InsertBefore("Z")
$markerRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$markerRange2 = $d.Range($insertPos, $insertPos + 1)
$markerRange2.Delete()
